$p = $ppt.ActivePresentation

# --- 1. Revert the "ERROR Device ... IO TABLE" slide addition: delete slide 3 ---
$p.Slides.Item(3).Delete()

# --- 2. Revert the date placeholder re-cache: 2023-11-16 -> 2023-11-02 ---
# The cached "datetimeFigureOut" text lives on the slide master and on every
# slide layout's date placeholder. Walk each one and fix the placeholder with
# PlaceholderFormat.Type = 16 (ppPlaceholderDate).
$newDate = "2023-11-02"
$ppPlaceholderDate = 16

$master = $p.Slides.Item(1).Master

function Update-DatePlaceholder($shapes, $newText, $dateType) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Type -eq 14) {
            if ($sh.PlaceholderFormat.Type -eq $dateType) {
                $sh.TextFrame.TextRange.Text = $newText
            }
        }
    }
}

Update-DatePlaceholder $master.Shapes $newDate $ppPlaceholderDate

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes $newDate $ppPlaceholderDate
}
